$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.008.00"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "'2.226.91"
$ws.Range("E3").Value = "  -0.58%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -1.46%  "

$ws.Range("D5").Value = "'299.18"
$ws.Range("E5").Value = "  -2.76%  "

$ws.Range("D6").Value = "'90.56"
$ws.Range("E6").Value = "  -3.97%  "

$ws.Range("D7").Value = "'0.557"
$ws.Range("E7").Value = "  -2.49%  "

$ws.Range("D8").Value = "'0.999"

$ws.Range("D9").Value = "'0.493"
$ws.Range("E9").Value = "  -5.55%  "

$ws.Range("D10").Value = "'33.08"
$ws.Range("E10").Value = "  -4.41%  "

$ws.Range("D11").Value = "'0.0778"
$ws.Range("E11").Value = "  -3.09%  "

$ws.Range("D12").Value = "'6.96"
$ws.Range("E12").Value = "  -3.42%  "

$ws.Range("E13").Value = "  -0.45%  "

$ws.Range("D14").Value = "'2.560.29"
$ws.Range("E14").Value = "  -0.82%  "

$ws.Range("D15").Value = "'2.219.08"
$ws.Range("E15").Value = "  -0.65%  "

$ws.Range("D16").Value = "'13.39"
$ws.Range("E16").Value = "  -1.05%  "

$ws.Range("D17").Value = "'0.777"
$ws.Range("E17").Value = "  -6.63%  "

$ws.Range("D18").Value = "'43.830.82"
$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").Value = "'11.79"
$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.0₃0905"
$ws.Range("E20").Value = "  -5.06%  "

$ws.Range("D21").Value = "'5.95"
$ws.Range("E21").Value = "  -6.30%  "

$ws.Range("D22").Value = "'64.59"
$ws.Range("E22").Value = "  -1.36%  "

$ws.Range("D23").Value = "'236.40"
$ws.Range("E23").Value = "  -0.54%  "

$ws.Range("D24").Value = "'2.82"
$ws.Range("E24").Value = "  -4.67%  "

$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("D26").Value = "'1.86"
$ws.Range("E26").Value = "  -5.58%  "

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'2.26"
$ws.Range("E27").Value = "  +1.70%  "

$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'38.73"
$ws.Range("E28").Value = "  +2.16%  "

$ws.Range("D29").Value = "'9.36"
$ws.Range("E29").Value = "  -4.28%  "

$ws.Range("D30").Value = "'152.66"
$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("D31").Value = "'19.19"
$ws.Range("E31").Value = "  -3.72%  "

$ws.Range("D32").Value = "'5.41"
$ws.Range("E32").Value = "  -8.90%  "

$ws.Range("D33").Value = "'0.0760"
$ws.Range("E33").Value = "  -4.22%  "

$ws.Range("D34").Value = "'2.48"
$ws.Range("E34").Value = "  -5.89%  "

$ws.Range("D35").Value = "'0.117"
$ws.Range("E35").Value = "  -1.69%  "

$ws.Range("D36").Value = "'2.85"
$ws.Range("E36").Value = "  -8.02%  "

$ws.Range("D37").Value = "'0.102"
$ws.Range("E37").Value = "  -6.71%  "

$ws.Range("D38").Value = "'1.69"
$ws.Range("E38").Value = "  -6.15%  "

$ws.Range("D39").Value = "'0.0300"
$ws.Range("E39").Value = "  +0.88%  "

$ws.Range("D40").Value = "'3.62"
$ws.Range("E40").Value = "  -3.57%  "

$ws.Range("D41").Value = "'3.16"
$ws.Range("E41").Value = "  -6.46%  "

$ws.Range("D42").Value = "'13.38"
$ws.Range("E42").Value = "  -10.21%  "

$ws.Range("E43").Value = "  -0.94%  "

$ws.Range("D44").Value = "'1.806.20"
$ws.Range("E44").Value = "  +1.34%  "

$ws.Range("D45").Value = "'1.80"
$ws.Range("E45").Value = "  +13.93%  "

$ws.Range("D46").Value = "'0.184"
$ws.Range("E46").Value = "  -4.04%  "

$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").Value = "'67.39"
$ws.Range("E47").Value = "  -3.86%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'94.44"
$ws.Range("E48").Value = "  -4.08%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'7.80"
$ws.Range("E49").Value = "  -3.27%  "

$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "'72.86"
$ws.Range("E50").Value = "  -7.52%  "

$ws.Range("D51").Value = "'4.60"
$ws.Range("E51").Value = "  -5.63%  "

